$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet's tab/sheet name
$ws.Name = "GammaFiber2F"

# Add a new row 16 with data, mirroring row 15's pattern (HexGrid-60degTilt5degRes)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
